# Auto-applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $r = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (prices, percents)
    # keep their exact original formatting (e.g. "587.55", "1.00.00", "0.999")
    # instead of being coerced into Excel numbers.
    $r.NumberFormat = "@"
    $r.Value = $Value
    # Drop back to the default style so no stray numFmt/style is left on the cell.
    $r.Style = "Normal"
}

Set-TextValue 'D2' '68.996.84'
Set-TextValue 'E2' '  -2.03%  '
Set-TextValue 'D3' '3.518.06'
Set-TextValue 'E3' '  -3.00%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '587.55'
Set-TextValue 'E5' '  +1.37%  '
Set-TextValue 'D6' '170.60'
Set-TextValue 'E6' '  -2.52%  '
Set-TextValue 'D7' '0.611'
Set-TextValue 'E7' '  +0.46%  '
Set-TextValue 'D8' '3.511.26'
Set-TextValue 'E8' '  -2.96%  '
Set-TextValue 'E9' '  +0.03%  '
Set-TextValue 'D10' '0.188'
Set-TextValue 'E10' '  -3.76%  '
Set-TextValue 'D11' '6.81'
Set-TextValue 'E11' '  -3.27%  '
Set-TextValue 'D12' '0.579'
Set-TextValue 'E12' '  -4.09%  '
Set-TextValue 'D13' '46.93'
Set-TextValue 'E13' '  -2.50%  '
Set-TextValue 'E14' '  -2.33%  '
Set-TextValue 'D15' '4.084.25'
Set-TextValue 'E15' '  -3.24%  '
Set-TextValue 'D16' '8.44'
Set-TextValue 'E16' '  -4.78%  '
Set-TextValue 'D17' '615.95'
Set-TextValue 'E17' '  -8.50%  '
Set-TextValue 'D18' '69.098.82'
Set-TextValue 'E18' '  -2.11%  '
Set-TextValue 'D19' '3.490.81'
Set-TextValue 'E19' '  -3.95%  '
Set-TextValue 'E20' '  -0.55%  '
Set-TextValue 'D21' '17.34'
Set-TextValue 'E21' '  -2.29%  '
Set-TextValue 'D22' '11.09'
Set-TextValue 'E22' '  -2.51%  '
Set-TextValue 'D23' '0.882'
Set-TextValue 'E23' '  -5.89%  '
Set-TextValue 'D24' '15.76'
Set-TextValue 'E24' '  -7.47%  '
Set-TextValue 'D25' '96.69'
Set-TextValue 'E25' '  -3.08%  '
Set-TextValue 'D26' '3.85'
Set-TextValue 'E26' '  -1.32%  '
Set-TextValue 'E27' '  +0.00%  '
Set-TextValue 'D28' '2.62'
Set-TextValue 'E28' '  -5.88%  '
Set-TextValue 'D29' '9.20'
Set-TextValue 'E29' '  -6.39%  '
Set-TextValue 'D30' '32.61'
Set-TextValue 'E30' '  -5.36%  '
Set-TextValue 'B31' 'Stacks'
Set-TextValue 'C31' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D31' '3.11'
Set-TextValue 'E31' '  -5.59%  '
Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '8.46'
Set-TextValue 'E32' '  -5.71%  '
Set-TextValue 'E33' '  -4.92%  '
Set-TextValue 'D34' '6.90'
Set-TextValue 'E34' '  -7.58%  '
Set-TextValue 'D35' '627.51'
Set-TextValue 'E35' '  +8.52%  '
Set-TextValue 'D36' '10.72'
Set-TextValue 'E36' '  -2.79%  '
Set-TextValue 'E37' '  -11.73%  '
Set-TextValue 'E38' '  -4.21%  '
Set-TextValue 'D39' '57.30'
Set-TextValue 'E39' '  -1.49%  '
Set-TextValue 'D40' '0.999'
Set-TextValue 'E40' '  -0.06%  '
Set-TextValue 'D41' '0.0449'
Set-TextValue 'E41' '  -0.56%  '
Set-TextValue 'E42' '  -4.10%  '
Set-TextValue 'D43' '3.368.06'
Set-TextValue 'D44' '0.326'
Set-TextValue 'E44' '  -4.88%  '
Set-TextValue 'D45' '32.65'
Set-TextValue 'E45' '  -5.20%  '
Set-TextValue 'D46' '0.0₃0693'
Set-TextValue 'E46' '  -5.00%  '
Set-TextValue 'E47' '  -5.21%  '
Set-TextValue 'E48' '  -1.70%  '
Set-TextValue 'E49' '  -2.78%  '
Set-TextValue 'D50' '133.13'
Set-TextValue 'E50' '  -2.37%  '
Set-TextValue 'E51' '  +12.50%  '
